$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")

# Insert a new row before row 29 (shifts rows 29-99 down to 30-100)
$ws.Rows("29:29").Insert()

# Set the new row's label cell
$ws.Range("A29").Value = "LOG prefix Name"

# Add the new comment describing the field
$comment = $ws.Range("A29").AddComment("Log prefix name to pass in the vsc config")
